$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 2).Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/major-diagnostic-category"
$meta.Cells.Item(3, 2).Value = "8.0.0"
$meta.Cells.Item(8, 2).Value = "2022-11-10T16:00:46+00:00"
$meta.Cells.Item(9, 2).Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Move the ele-1/ext-1 constraint text off the "Extension" summary row (row 2)
# and onto the "Extension.extension" row (row 4), where it belongs.
$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Cells.Item(2, 35).Value = ""
$elements.Cells.Item(4, 35).Value = $constraintText

# Update the base URL referenced as the Fixed Value on Extension.url
$elements.Cells.Item(5, 17).Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/major-diagnostic-category"

# Update the Binding Value Set URL/name on the valueCodeableConcept slice
$elements.Cells.Item(7, 25).Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/payer-major-diagnostic-category"

# Widen column Y (Binding Value Set) to fit the new URL text
$elements.Columns.Item(25).ColumnWidth = 70
